# create new Tweet functionality with their hashtags saved with tweet
#
# This script mutates tasks.xlsx (Sheet1) to:
#  - mark "Give follow functionality" / "Give unfollow functionality" (row 54/55)
#    and "Give Tweet Adding Functionality" (row 61) as done (red -> normal/black),
#    by copying the formatting already used by "done" rows.
#  - insert a new blank separator row (row 60) using the same formatting as the
#    other blank separator rows (56/58).
#  - append a new block of "to do" (red) tasks about tweets (rows 71-75), reusing
#    the formatting already used for the existing red "to do" rows.
#  - append three new plain follow-up notes (rows 77, 79, 81).
#  - update the active selection to D63, matching where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark existing follow/unfollow tasks as completed (row 54 & 55: B:D red -> normal) ---
$ws.Range("E54:F54").Copy($ws.Range("B54:D54"))
$ws.Range("E55:F55").Copy($ws.Range("B55:D55"))

# --- New blank separator row 60 (same look as the other blank rows, e.g. row 56) ---
$ws.Range("B56:D56").Copy($ws.Range("B60:D60"))

# --- Mark "Give Tweet Adding Functionality" as completed (row 61: B:D red -> normal) ---
$ws.Range("E54:F54").Copy($ws.Range("B61:D61"))

# --- New "to do" block about tweets (rows 71-75), styled like the existing red block ---
$ws.Range("B53:F53").Copy($ws.Range("B71:F71"))
$ws.Range("B71").Value = "Make tweets in reverse chronological order"

$ws.Range("B53:F53").Copy($ws.Range("B72:F72"))

$ws.Range("B53:F53").Copy($ws.Range("B73:F73"))
$ws.Range("B73").Value = "Like/Dislike tweet - reaction of tweet"

$ws.Range("B53:F53").Copy($ws.Range("B74:F74"))

$ws.Range("B53:F53").Copy($ws.Range("B75:F75"))
$ws.Range("B75").Value = "Edit/Delete tweet option to the tweet owner"

# --- New plain follow-up notes ---
$ws.Range("B77").Value = "Tweet Hashtags and Tweet User handling"
$ws.Range("B79").Value = "Duplcate entries in follower table handling put index on entries"
$ws.Range("B81").Value = "Hashtags are case insensitive - I think its good"

# --- Update selection to where the author left off ---
$ws.Range("D63").Select()
